# Refresh the scraped "Price" (D) and "Volume(1h)" (E) columns for every
# coin row (2-51) with the latest values from the source feed.
#
# Price cells are stored as literal TEXT in the source data (not numbers),
# e.g. "309.66" or "26.662.37" (a thousands-grouped integer price). Excel
# auto-converts plain numeric-looking input to a Number, so values that
# would parse as a float get a leading apostrophe to force text entry -
# this preserves exact formatting (trailing zeros, decimal precision)
# and keeps the cell type as text, matching the source.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.662.37'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '1.827.56'
$ws.Range("E3").Value = '  +1.82%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '''309.66'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").Value = '''1.008'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '''0.4663'
$ws.Range("E7").Value = '  +3.36%  '
$ws.Range("D8").Value = '''0.3598'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.07143'
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").Value = '''0.9050'
$ws.Range("E10").Value = '  +2.41%  '
$ws.Range("D11").Value = '''0.07702'
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '1.860.19'
$ws.Range("E13").Value = '  +5.15%  '
$ws.Range("D14").Value = '''5.270'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '''6.377'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '''87.77'
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").Value = '26.698.24'
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '''1.912'
$ws.Range("E24").Value = '  -2.86%  '
$ws.Range("D25").Value = '''152.93'
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("E26").Value = '  +0.73%  '
$ws.Range("D27").Value = '''1.997'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").Value = '''113.84'
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").Value = '''4.873'
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").Value = '''0.08819'
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("D31").Value = '''3.124'
$ws.Range("D32").Value = '''2.843'
$ws.Range("E32").Value = '  +3.22%  '
$ws.Range("D33").Value = '''1.166'
$ws.Range("E33").Value = '  +5.67%  '
$ws.Range("D34").Value = '''0.7361'
$ws.Range("E34").Value = '  +1.75%  '
$ws.Range("D35").Value = '''4.440'
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").Value = '''1.082'
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("D37").Value = '''0.01935'
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").Value = '''0.05161'
$ws.Range("E38").Value = '  +1.29%  '
$ws.Range("D39").Value = '''2.924'
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("D40").Value = '''6.878'
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("D41").Value = '''0.5065'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").Value = '''0.4668'
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("D46").Value = '''10.05'
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").Value = '''98.60'
$ws.Range("E47").Value = '  -2.31%  '
$ws.Range("D48").Value = '''1.575'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").Value = '''0.06035'
$ws.Range("E49").Value = '  +1.23%  '
$ws.Range("D50").Value = '''63.95'
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("D51").Value = '''35.80'
$ws.Range("E51").Value = '  -0.69%  '
